$d = $word.ActiveDocument

# -------------------------------------------------------------------
# 1) Remove the "_GoBack" bookmark from the title paragraph.
#    (Bookmarks.Add() below will naturally re-create it at the new
#    location; deleting it here also causes the remaining
#    "_abifqbd3ehrc" bookmark to be renumbered from id=2 to id=1,
#    matching the target XML.)
# -------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# -------------------------------------------------------------------
# 2) Split the N02 bullet's run into "N02: " + the rest, and update
#    the wording from "Criação do banco de dados dos clientes." to
#    "Criação do cadastro para os clientes."
# -------------------------------------------------------------------
$d.Content.Find.Execute(
    "N02: Criação do banco de dados dos clientes.", $true, $false, $false, $false, $false,
    $true, 1, $false, "N02: Criação do cadastro para os clientes.", 2) | Out-Null

# Locate the N02 paragraph again (text has just changed) and split the
# single run into two runs right after "N02: " (5 characters) by
# inserting a paragraph break there and immediately merging the break
# back out. That leaves two clean <w:r> runs with no leftover explicit
# run formatting, while preserving the paragraph's own identity.
$n02Para = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "N02:*") {
        $n02Para = $para
        break
    }
}
$n02Start = $n02Para.Range.Start
$splitPoint = $d.Range($n02Start + 5, $n02Start + 5)
$splitPoint.InsertParagraphAfter()

# Re-fetch the (now two) paragraphs and merge them back into one by
# deleting the paragraph mark that separates them - this keeps the
# original paragraph (with its numbering / style) but now holding two
# separate runs: "N02: " and "Criação do cadastro para os clientes."
$n02Para = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "N02:*") {
        $n02Para = $para
        break
    }
}
$markEnd = $n02Para.Range.End
$markRange = $d.Range($markEnd - 1, $markEnd)
$markRange.Delete()

# -------------------------------------------------------------------
# 3) Add a fresh "_GoBack" bookmark at the very start of the N03
#    bullet (collapsed / zero-length), matching the relocated marker
#    in the target document.
# -------------------------------------------------------------------
$n03Para = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "N03:*") {
        $n03Para = $para
        break
    }
}
$n03Start = $n03Para.Range.Start
$n03Collapsed = $d.Range($n03Start, $n03Start)
$d.Bookmarks.Add("_GoBack", $n03Collapsed) | Out-Null
